$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1223.3572
$ws.Range("I19").Value = 1571.4
$ws.Range("J19").Value = 353.25
$ws.Range("K19").Value = 1571.4
$ws.Range("L19").Value = 353.25
$ws.Range("M19").Value = -1396.4
$ws.Range("N19").Value = -703.25

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H43").Value = 1911.8572
$ws.Range("J43").Value = 1949
$ws.Range("L43").Value = 1949
$ws.Range("N43").Value = -2087

$ws.Range("H98").Value = 13243.25
$ws.Range("I98").Value = 10611.5
$ws.Range("K98").Value = 10611.5
$ws.Range("M98").Value = -9113.5

$ws.Range("H122").Value = 13243.25
$ws.Range("I122").Value = 10611.5
$ws.Range("K122").Value = 31834.5
$ws.Range("M122").Value = -29384.5

$ws.Range("H129").Value = 2942911.2
$ws.Range("J129").Value = 1821.25
$ws.Range("L129").Value = 5463.75
$ws.Range("N129").Value = -15463.75

$ws.Range("H134").Value = 77000
$ws.Range("J134").Value = 77000
$ws.Range("L134").Value = 77000
$ws.Range("N134").Value = -87140

$ws.Range("H138").Value = 7695352.5
$ws.Range("J138").Value = 2948.125
$ws.Range("L138").Value = 8844.375
$ws.Range("N138").Value = -19124.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 3000
$ws.Range("I41").Value = 3000
$ws.Range("K41").Value = 3000
$ws.Range("M41").Value = -2586

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H96").Value = 27500
$ws.Range("J96").Value = 27500
$ws.Range("L96").Value = 27500
$ws.Range("N96").Value = -32992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8066.3335
$ws.Range("I134").Value = 6100
$ws.Range("J134").Value = 11999
$ws.Range("K134").Value = 18300
$ws.Range("L134").Value = 35997
$ws.Range("M134").Value = -15765
$ws.Range("N134").Value = -41067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1196.8
$ws.Range("J2").Value = 1811.3334
$ws.Range("L2").Value = 10868.0004
$ws.Range("N2").Value = -11094.0004

$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 5000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -14888
$ws.Range("N3").ClearContents()

$ws.Range("H17").Value = 429.75
$ws.Range("I17").Value = 358
$ws.Range("K17").Value = 1074
$ws.Range("M17").Value = -905

$ws.Range("H34").Value = 1000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3000
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -3168

$ws.Range("H46").Value = 519
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H81").Value = 9000
$ws.Range("I81").Value = 9000
$ws.Range("K81").Value = 27000
$ws.Range("M81").Value = -25877

$ws.Range("H84").Value = 9000
$ws.Range("I84").Value = 9000
$ws.Range("K84").Value = 81000
$ws.Range("M84").Value = -75384

$ws.Range("H103").Value = 2600
$ws.Range("I103").Value = 200
$ws.Range("J103").Value = 5000
$ws.Range("K103").Value = 600
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = 279
$ws.Range("N103").Value = -16758

$ws.Range("H104").Value = 5999.6665
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 5999.6665
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 17998.9995
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -23240.9995

$ws.Range("H114").Value = 1513.3334
$ws.Range("I114").Value = 1513.3334
$ws.Range("K114").Value = 4540.0002
$ws.Range("M114").Value = -1286.0002

$ws.Range("H131").Value = 935
$ws.Range("I131").Value = 935
$ws.Range("K131").Value = 2805
$ws.Range("M131").Value = 2235

$ws.Range("H133").Value = 1000
$ws.Range("I133").Value = 1000
$ws.Range("K133").Value = 3000
$ws.Range("M133").Value = 2060

$ws.Range("H138").Value = 3445
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3445
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 10335
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -20615

$ws.Range("H140").Value = 2074.1428
$ws.Range("I140").Value = 2074.1428
$ws.Range("K140").Value = 6222.428400000001
$ws.Range("M140").Value = -1042.428400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 20000
$ws.Range("J48").Value = 20000
$ws.Range("L48").Value = 20000
$ws.Range("N48").Value = -20970

$ws.Range("H122").Value = 999.8
$ws.Range("I122").Value = 999.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2999.4
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -549.3999999999996
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 3794
$ws.Range("I132").Value = 3794
$ws.Range("K132").Value = 11382
$ws.Range("M132").Value = -8852

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H46").Value = 7900
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 20000
$ws.Range("J70").Value = 20000
$ws.Range("L70").Value = 20000
$ws.Range("N70").Value = -20630

$ws.Range("H73").Value = 20000
$ws.Range("J73").Value = 20000
$ws.Range("L73").Value = 20000
$ws.Range("N73").Value = -22184

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

